$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "alvaro"

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "amigo"

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "amiga"

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "amiga"

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "amiga"

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "amiga"

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "amiga"

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "amiga"

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "amiga"

$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "amiga"

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "amiga"

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "amiga"

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "amiga"

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "amiga"

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "amiga"

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "amiga"

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "amigas"

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "amigas2"

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "amigos2"

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "ala"

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "aladin"

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "genio"

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "genios"

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "genioso"
$ws.Cells.Item(41, 4).Value = 999999999

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "abuela"
$ws.Cells.Item(42, 4).Value = 333333333

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "abuelo"
$ws.Cells.Item(43, 4).Value = 333333333

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "abueloa"
$ws.Cells.Item(44, 4).Value = 333333333

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "abuelos"
$ws.Cells.Item(45, 4).Value = 333333333

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "pere"
$ws.Cells.Item(46, 4).Value = 333333333

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "joan"
